$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values currently stored in A1:A10 (top to bottom) before
# clearing them, so the same shared-string content can be re-written
# into the new layout.
$values = @(
    $ws.Range("A1").Value(),
    $ws.Range("A2").Value(),
    $ws.Range("A3").Value(),
    $ws.Range("A4").Value(),
    $ws.Range("A5").Value(),
    $ws.Range("A6").Value(),
    $ws.Range("A7").Value(),
    $ws.Range("A8").Value(),
    $ws.Range("A9").Value(),
    $ws.Range("A10").Value()
)

# Remove the old single-column layout in A1:A10.
$ws.Range("A1:A10").Clear()

# Re-write the same values across row 1, columns C through L, preserving
# the wrap-text formatting that the cells had before.
$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L")
for ($i = 0; $i -lt $values.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "1")
    $cell.Value = $values[$i]
    $cell.WrapText = $true
}

# Restore automatic row height for row 1 (writing the long wrapped text
# forces an explicit custom height; AutoFit puts it back to the default).
$ws.Rows("1").AutoFit()
